# Update load_dataset and add use-case datasets
#
# The table row that described the "t1_overview_proteomics" page is renamed
# to "t1_overview_datasets", and its description is broadened from
# "Proteomic example datasets" to "Omics example datasets".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nameCell = $ws.Cells.Find("t1_overview_proteomics")
if ($nameCell -ne $null) {
    $nameCell.Value = "t1_overview_datasets"
} else {
    $ws.Range("A2").Value = "t1_overview_datasets"
}

$descCell = $ws.Cells.Find("Proteomic example datasets")
if ($descCell -ne $null) {
    $descCell.Value = "Omics example datasets"
} else {
    $ws.Range("B2").Value = "Omics example datasets"
}

# Move the active cell selection from B8 to B5
$ws.Range("B5").Select()
